$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.76479233333333
$ws.Range("H2").Value = 56.294377
$ws.Range("I2").Value = 0.09818846546758657
$ws.Range("J2").Value = 0.09818846546758656
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.001642333333333
$ws.Range("N2").Value = 9.004927
$ws.Range("O2").Value = 0.05169795991651582
$ws.Range("P2").Value = 0.05169795991651582
$ws.Range("Q2").Value = 56.32519504394211
$ws.Range("R2").Value = 506.926755395479
$ws.Range("S2").Value = 0.005076143352007488
$ws.Range("T2").Value = 0.005076143352007488

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.76479233333333
$ws.Range("H3").Value = 56.294377
$ws.Range("I3").Value = 0.09818846546758657
$ws.Range("J3").Value = 0.09818846546758656
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 11.05428166666667
$ws.Range("N3").Value = 33.162845
$ws.Range("O3").Value = 0.1903903753498087
$ws.Range("P3").Value = 0.1903903753498088
$ws.Range("Q3").Value = 207.4312998691739
$ws.Range("R3").Value = 1866.881698822565
$ws.Range("S3").Value = 0.01869413879539554
$ws.Range("T3").Value = 0.01869413879539554

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.76479233333333
$ws.Range("H4").Value = 56.294377
$ws.Range("I4").Value = 0.09818846546758657
$ws.Range("J4").Value = 0.09818846546758656
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.49405866666666
$ws.Range("N4").Value = 49.482176
$ws.Range("O4").Value = 0.2840808761059341
$ws.Range("P4").Value = 0.2840808761059341
$ws.Range("Q4").Value = 309.5075856138168
$ws.Range("R4").Value = 2785.568270524352
$ws.Range("S4").Value = 0.02789346529352925
$ws.Range("T4").Value = 0.02789346529352925

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.76479233333333
$ws.Range("H5").Value = 56.294377
$ws.Range("I5").Value = 0.09818846546758657
$ws.Range("J5").Value = 0.09818846546758656
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 27.51115433333333
$ws.Range("N5").Value = 82.533463
$ws.Range("O5").Value = 0.4738307886277414
$ws.Range("P5").Value = 0.4738307886277414
$ws.Range("Q5").Value = 516.2410979152834
$ws.Range("R5").Value = 4646.169881237551
$ws.Range("S5").Value = 0.0465247180266543
$ws.Range("T5").Value = 0.04652471802665429

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 72.11798333333333
$ws.Range("H6").Value = 216.35395
$ws.Range("I6").Value = 0.3773638413007209
$ws.Range("J6").Value = 0.3773638413007209
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.001642333333333
$ws.Range("N6").Value = 9.004927
$ws.Range("O6").Value = 0.05169795991651582
$ws.Range("P6").Value = 0.05169795991651582
$ws.Range("Q6").Value = 216.4723917679611
$ws.Range("R6").Value = 1948.25152591165
$ws.Range("S6").Value = 0.0195089407415071
$ws.Range("T6").Value = 0.0195089407415071

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 72.11798333333333
$ws.Range("H7").Value = 216.35395
$ws.Range("I7").Value = 0.3773638413007209
$ws.Range("J7").Value = 0.3773638413007209
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 11.05428166666667
$ws.Range("N7").Value = 33.162845
$ws.Range("O7").Value = 0.1903903753498087
$ws.Range("P7").Value = 0.1903903753498088
$ws.Range("Q7").Value = 797.2125009986389
$ws.Range("R7").Value = 7174.912508987751
$ws.Range("S7").Value = 0.07184644338868991
$ws.Range("T7").Value = 0.07184644338868991

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 72.11798333333333
$ws.Range("H8").Value = 216.35395
$ws.Range("I8").Value = 0.3773638413007209
$ws.Range("J8").Value = 0.3773638413007209
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.49405866666666
$ws.Range("N8").Value = 49.482176
$ws.Range("O8").Value = 0.2840808761059341
$ws.Range("P8").Value = 0.2840808761059341
$ws.Range("Q8").Value = 1189.518248021689
$ws.Range("R8").Value = 10705.6642321952
$ws.Range("S8").Value = 0.1072018506474095
$ws.Range("T8").Value = 0.1072018506474095

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 72.11798333333333
$ws.Range("H9").Value = 216.35395
$ws.Range("I9").Value = 0.3773638413007209
$ws.Range("J9").Value = 0.3773638413007209
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.51115433333333
$ws.Range("N9").Value = 82.533463
$ws.Range("O9").Value = 0.4738307886277414
$ws.Range("P9").Value = 0.4738307886277414
$ws.Range("Q9").Value = 1984.048969692094
$ws.Range("R9").Value = 17856.44072722885
$ws.Range("S9").Value = 0.1788066065231144
$ws.Range("T9").Value = 0.1788066065231144

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 83.31930033333333
$ws.Range("H10").Value = 249.957901
$ws.Range("I10").Value = 0.4359757410707099
$ws.Range("J10").Value = 0.4359757410707098
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.001642333333333
$ws.Range("N10").Value = 9.004927
$ws.Range("O10").Value = 0.05169795991651582
$ws.Range("P10").Value = 0.05169795991651582
$ws.Range("Q10").Value = 250.0947390642474
$ws.Range("R10").Value = 2250.852651578227
$ws.Range("S10").Value = 0.02253905638644684
$ws.Range("T10").Value = 0.02253905638644684

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 83.31930033333333
$ws.Range("H11").Value = 249.957901
$ws.Range("I11").Value = 0.4359757410707099
$ws.Range("J11").Value = 0.4359757410707098
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 11.05428166666667
$ws.Range("N11").Value = 33.162845
$ws.Range("O11").Value = 0.1903903753498087
$ws.Range("P11").Value = 0.1903903753498088
$ws.Range("Q11").Value = 921.0350141542607
$ws.Range("R11").Value = 8289.315127388345
$ws.Range("S11").Value = 0.08300558498586348
$ws.Range("T11").Value = 0.08300558498586348

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 83.31930033333333
$ws.Range("H12").Value = 249.957901
$ws.Range("I12").Value = 0.4359757410707099
$ws.Range("J12").Value = 0.4359757410707098
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 16.49405866666666
$ws.Range("N12").Value = 49.482176
$ws.Range("O12").Value = 0.2840808761059341
$ws.Range("P12").Value = 0.2840808761059341
$ws.Range("Q12").Value = 1374.273427763619
$ws.Range("R12").Value = 12368.46084987257
$ws.Range("S12").Value = 0.1238523704843011
$ws.Range("T12").Value = 0.1238523704843011

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 83.31930033333333
$ws.Range("H13").Value = 249.957901
$ws.Range("I13").Value = 0.4359757410707099
$ws.Range("J13").Value = 0.4359757410707098
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 27.51115433333333
$ws.Range("N13").Value = 82.533463
$ws.Range("O13").Value = 0.4738307886277414
$ws.Range("P13").Value = 0.4738307886277414
$ws.Range("Q13").Value = 2292.210130415685
$ws.Range("R13").Value = 20629.89117374116
$ws.Range("S13").Value = 0.2065787292140984
$ws.Range("T13").Value = 0.2065787292140984

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 16.90787
$ws.Range("H14").Value = 50.72361
$ws.Range("I14").Value = 0.08847195216098278
$ws.Range("J14").Value = 0.08847195216098278
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.001642333333333
$ws.Range("N14").Value = 9.004927
$ws.Range("O14").Value = 0.05169795991651582
$ws.Range("P14").Value = 0.05169795991651582
$ws.Range("Q14").Value = 50.75137835849667
$ws.Range("R14").Value = 456.76240522647
$ws.Range("S14").Value = 0.004573819436554392
$ws.Range("T14").Value = 0.004573819436554393

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 16.90787
$ws.Range("H15").Value = 50.72361
$ws.Range("I15").Value = 0.08847195216098278
$ws.Range("J15").Value = 0.08847195216098278
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 11.05428166666667
$ws.Range("N15").Value = 33.162845
$ws.Range("O15").Value = 0.1903903753498087
$ws.Range("P15").Value = 0.1903903753498088
$ws.Range("Q15").Value = 186.9043573633834
$ws.Range("R15").Value = 1682.13921627045
$ws.Range("S15").Value = 0.01684420817985983
$ws.Range("T15").Value = 0.01684420817985983

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 16.90787
$ws.Range("H16").Value = 50.72361
$ws.Range("I16").Value = 0.08847195216098278
$ws.Range("J16").Value = 0.08847195216098278
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 16.49405866666666
$ws.Range("N16").Value = 49.482176
$ws.Range("O16").Value = 0.2840808761059341
$ws.Range("P16").Value = 0.2840808761059341
$ws.Range("Q16").Value = 278.8793997083733
$ws.Range("R16").Value = 2509.91459737536
$ws.Range("S16").Value = 0.02513318968069428
$ws.Range("T16").Value = 0.02513318968069428

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 16.90787
$ws.Range("H17").Value = 50.72361
$ws.Range("I17").Value = 0.08847195216098278
$ws.Range("J17").Value = 0.08847195216098278
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 27.51115433333333
$ws.Range("N17").Value = 82.533463
$ws.Range("O17").Value = 0.4738307886277414
$ws.Range("P17").Value = 0.4738307886277414
$ws.Range("Q17").Value = 465.1550210179366
$ws.Range("R17").Value = 4186.39518916143
$ws.Range("S17").Value = 0.04192073486387427
$ws.Range("T17").Value = 0.04192073486387428
